$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.155.22'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.34%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.880.22'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.50%  '

# Row 4
$ws.Range("E4").Value = '  +0.16%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.27%  '

# Row 6
$ws.Range("E6").Value = '  +0.12%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5083'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.32%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3857'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.82%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09125'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.38%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.124'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.54%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.57'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.66%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.356'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.72%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.79'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.47%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.873.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.97%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.207'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.44%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.002'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.18%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001113'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.00%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.25'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.66%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06612'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.13%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.41%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.21%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.118'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.31%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.189.38'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.39%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.45'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.22%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.280'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.71%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.572'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.79%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.091.45'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.89%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.81'
$ws.Range("D28").Style = "Normal"

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '157.30'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.30%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.74'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.49%  '

# Row 31
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.066'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.30%  '

# Row 32
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1054'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.71%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.619'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.40%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.603'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.43%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.716'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.23%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02475'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.44%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06579'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.32%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2179'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.49%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.215'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.85%  '

# Row 40
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6419'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.02%  '

# Row 41
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.58'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.76%  '

# Row 42
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.239'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.84%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.935'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.56%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.26'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.56%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6013'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.10%  '

# Row 46
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.677'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.14%  '

# Row 47
$ws.Range("B47").Value = 'WEMIXTOKEN'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.275'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.09%  '

# Row 48
$ws.Range("B48").Value = 'EOS'
$ws.Range("C48").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.233'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.93%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.004'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.80%  '

# Row 50
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '121.46'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.20%  '

# Row 51
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.94'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.98%  '
